# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K"), rows 2-11
$newK = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 3
    6  = 1
    7  = 1
    8  = 1
    9  = 2
    10 = 1
    11 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
